$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 239, shifting rows 239:341 down to 240:342
$ws.Rows.Item(239).Insert()

# Fill the new row 239 with the new data row
$ws.Cells.Item(239, 1).Value = 3
$ws.Cells.Item(239, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(239, 3).Value = "Coquimbo"
$ws.Cells.Item(239, 4).Value = 44704
$ws.Cells.Item(239, 5).Value = 5
$ws.Cells.Item(239, 6).Value = 100112012
$ws.Cells.Item(239, 7).Value = "Espinaca"
$ws.Cells.Item(239, 8).Value = "Sin especificar"
$ws.Cells.Item(239, 9).Value = "Primera"
$ws.Cells.Item(239, 10).Value = 230
$ws.Cells.Item(239, 11).Value = 3500
$ws.Cells.Item(239, 12).Value = 3800
$ws.Cells.Item(239, 13).Value = 3657
$ws.Cells.Item(239, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(239, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(239, 16).Value = 1219
$ws.Cells.Item(239, 17).Value = 3
$ws.Cells.Item(239, 18).Value = "Hortaliza"
